# Arregle un poco el formato de las tablas y agregue la variable 'partido lider'

$wb = $excel.ActiveWorkbook
$wsBase = $wb.Worksheets.Item("Base de datos")
$wsElec = $wb.Worksheets.Item("Elecciones")

# --- Sheet "Elecciones" ----------------------------------------------------
# Pre-seed the text that will later live in L2, writing it into the
# still-unshifted K2 cell so the shared-string table keeps the same layout
# the workbook ends up with once everything below is applied.
$wsElec.Range("K2").Value = "Partido Lider"

# Insert a new column at H ("Partido Lider"), pushing siglas (H->I),
# Numero de votos (I->J), arrastre (J->K), cantidad de colegios (K->L),
# cantidad inscritos (L->M), votos validos (M->N) and votos nulos (N->O)
# one column to the right.
$wsElec.Columns.Item(8).Insert()

$wsElec.Range("H1").Value = "Partido Lider"
$wsElec.Range("H2").Value = "PLD"
$wsElec.Range("H3").Value = "PLD"
$wsElec.Range("H4").Value = "PLD"
$wsElec.Range("H5").Value = "PLD"
$wsElec.Range("H6").Value = "PLD"
$wsElec.Range("H7").Value = "PLD"

# Fix the text that got displaced into L2 by the column insert above.
$wsElec.Range("L2").Value = "Valido solo desde 1998"

# Rename headers that shifted right.
$wsElec.Range("J1").Value = "# de votos"
$wsElec.Range("L1").Value = "# de colegios"
$wsElec.Range("M1").Value = "# inscritos"

# Best-effort column widths (auto-fit look) for the shifted columns.
$wsElec.Columns.Item(9).ColumnWidth = 14.5
$wsElec.Columns.Item(10).ColumnWidth = 14.5
$wsElec.Columns.Item(11).ColumnWidth = 6.6666666666666667
$wsElec.Columns.Item(12).ColumnWidth = 19.166666666666668
$wsElec.Columns.Item(13).ColumnWidth = 10.833333333333334
$wsElec.Columns.Item(14).ColumnWidth = 9.5

# --- Sheet "Base de datos": reset the view (drop topLeftCell/selection) ---
$wsBase.Activate()
$wsBase.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# --- Leave "Elecciones" as the active/visible tab, selection at O1 --------
$wsElec.Activate()
$wsElec.Range("O1").Select()
